# market_health_data.xlsx update: 2025-11-05 14:24
# - Metadata: bump "Last Updated" timestamp by a minute
# - Industry Analysis: refresh the "1 Year" (column F) figures for rows 2-76
# - Stock List: a new ticker (CAPTRU-RE1) enters at the top of the list,
#   pushing every other row down by one and dropping the last row (TRAVELFOOD)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata!A2 - "Last Updated" timestamp
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 02:23 PM"

# ---------------------------------------------------------------------------
# 2. Industry Analysis - column F ("1 Year") refreshed values, rows 2-76
# ---------------------------------------------------------------------------
$industry = $wb.Worksheets.Item("Industry Analysis")

$industryF = @{
    2  = 21.3
    3  = -4.3927
    4  = 35.9445
    5  = -51.0482
    6  = 57.2275
    7  = -9.640700000000001
    8  = -6.1449
    9  = 36.9733
    10 = -4.7026
    11 = 46.5317
    12 = -2.102
    13 = 17.4681
    14 = -33.0245
    15 = 1.0205
    16 = 2.0426
    17 = -16.2411
    18 = 7.4627
    19 = -25.798
    20 = 47.7485
    21 = 19.5587
    22 = 76.5603
    23 = -54.2675
    24 = -0.8811
    25 = 4.8518
    26 = 3.6831
    27 = -34.0874
    28 = -11.9893
    29 = -12.994
    30 = 25.5415
    31 = 56.5088
    32 = 2.0908
    33 = -4.7193
    34 = 22.8807
    35 = 5.3359
    36 = -5.1995
    37 = -5.6238
    38 = -22.595
    39 = 10.8405
    40 = -7.5963
    41 = -4.552
    42 = 22.3098
    43 = 14.0694
    44 = -9.6066
    45 = 27.639
    46 = -6.3484
    47 = -40.5302
    48 = -29.7988
    49 = -24.0791
    50 = -49.1803
    51 = -51.6023
    52 = -34.4756
    53 = -11.5478
    54 = -2.3796
    55 = -15.4382
    56 = -27.6987
    57 = -27.1559
    58 = -2.1585
    59 = -23.0964
    60 = -13.3217
    61 = -8.1496
    62 = -16.0695
    63 = -12.5465
    64 = 47.7264
    65 = -42.4232
    66 = 11.3291
    67 = 14.3746
    68 = 32.6702
    69 = -17.0097
    70 = -13.5162
    71 = 11.4259
    72 = 2.6754
    73 = -11.1574
    74 = -13.2502
    75 = 24.7078
    76 = 53.3554
}

foreach ($row in $industryF.Keys) {
    $industry.Cells.Item($row, 6).Value = $industryF[$row]
}

# ---------------------------------------------------------------------------
# 3. Stock List - insert new row 2 (CAPTRU-RE1), shift everything else down
#    one row, and drop the old last row (previously TRAVELFOOD at row 76).
# ---------------------------------------------------------------------------
$stocks = $wb.Worksheets.Item("Stock List")

$stocks.Rows.Item(2).Insert()
$stocks.Rows.Item(2).ClearFormats()

$stocks.Range("A2").Value = [char]0x1F4CB
$stocks.Range("B2").Value = "CAPTRU-RE1"
$stocks.Range("C2").Value = "CAPTRU-RE1"
$stocks.Range("D2").Value = 5.67
$stocks.Range("E2").Value = -11.9565
$stocks.Range("F2").Value = "N/A"
$stocks.Range("G2").Value = "N/A"
$stocks.Range("H2").Value = 0

# the insert pushed the former last row (TRAVELFOOD) down to row 77; remove it
# so the sheet stays at 76 rows (header + 75 data rows).
$stocks.Rows.Item(77).Delete()
